# Update "Forecast Comparison" data: the weekly forecast table rolled forward
# by one week (new week added at the bottom, oldest week dropped) and the
# forecast numbers were refreshed. The "Summary" sheet metrics were
# recalculated to match.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($range, [string]$text)
    # Force the cell to stay text (avoids Excel auto-converting values that
    # look like dates/numbers into real dates/numbers), while preserving the
    # cell's original style.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

# ---------------------------------------------------------------------------
# Sheet 1: Forecast Comparison
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$forecastRows = @(
    @{ Row = 2;  Date = "2025-02-02"; D = 41; E = 42; F = 50; G = 58; H = 69 },
    @{ Row = 3;  Date = "2025-02-09"; D = 34; E = 36; F = 43; G = 50; H = 61 },
    @{ Row = 4;  Date = "2025-02-16"; D = 35; E = 36; F = 43; G = 50; H = 61 },
    @{ Row = 5;  Date = "2025-02-23"; D = 47; E = 38; F = 46; G = 53; H = 65 },
    @{ Row = 6;  Date = "2025-03-02"; D = 49; E = 38; F = 46; G = 54; H = 66 },
    @{ Row = 7;  Date = "2025-03-09"; D = 48; E = 37; F = 45; G = 53; H = 65 },
    @{ Row = 8;  Date = "2025-03-16"; D = 51; E = 39; F = 47; G = 55; H = 69 },
    @{ Row = 9;  Date = "2025-03-23"; D = 51; E = 39; F = 48; G = 58; H = 75 },
    @{ Row = 10; Date = "2025-03-30"; D = 49; E = 38; F = 45; G = 54; H = 67 },
    @{ Row = 11; Date = "2025-04-06"; D = 49; E = 38; F = 46; G = 56; H = 72 },
    @{ Row = 12; Date = "2025-04-13"; D = 48; E = 37; F = 45; G = 54; H = 69 },
    @{ Row = 13; Date = "2025-04-20"; D = 49; E = 38; F = 46; G = 56; H = 72 },
    @{ Row = 14; Date = "2025-04-27"; D = 47; E = 36; F = 44; G = 54; H = 69 },
    @{ Row = 15; Date = "2025-05-04"; D = 44; E = 36; F = 44; G = 54; H = 70 },
    @{ Row = 16; Date = "2025-05-11"; D = 40; E = 36; F = 44; G = 54; H = 70 },
    @{ Row = 17; Date = "2025-05-18"; D = 38; E = 35; F = 43; G = 53; H = 70 }
)

foreach ($r in $forecastRows) {
    Set-TextCell $ws1.Range("B$($r.Row)") $r.Date
    $ws1.Range("D$($r.Row)").Value = $r.D
    $ws1.Range("E$($r.Row)").Value = $r.E
    $ws1.Range("F$($r.Row)").Value = $r.F
    $ws1.Range("G$($r.Row)").Value = $r.G
    $ws1.Range("H$($r.Row)").Value = $r.H
}

# ---------------------------------------------------------------------------
# Sheet 2: Summary
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

Set-TextCell $ws2.Range("B2") "2022-12-25 to 2025-01-26"
Set-TextCell $ws2.Range("B4") "277"
Set-TextCell $ws2.Range("B6") "89"
Set-TextCell $ws2.Range("B8") "10823 units"
Set-TextCell $ws2.Range("B9") "721"
Set-TextCell $ws2.Range("B10") "356"
Set-TextCell $ws2.Range("B11") "157"
Set-TextCell $ws2.Range("B12") "51"
Set-TextCell $ws2.Range("B14") "34"
